# Update the "Last Updated" timestamp on the Metadata sheet
$wb = $excel.ActiveWorkbook
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Cells.Item(2, 1).Value = "29 Oct 2025, 07:17 PM"

# Insert a new "Top Losers" worksheet right after "Top Gainers"
$gainersWs = $wb.Worksheets.Item("Top Gainers")
$losersWs = $wb.Worksheets.Add($null, $gainersWs)
$losersWs.Name = "Top Losers"

# Header row
$losersWs.Cells.Item(1, 1).Value = "Icon"
$losersWs.Cells.Item(1, 2).Value = "Stock"
$losersWs.Cells.Item(1, 3).Value = "Latest"
$losersWs.Cells.Item(1, 4).Value = "Weekly"
$losersWs.Cells.Item(1, 5).Value = "Monthly"

$headerRange = $losersWs.Range("A1:E1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows: Icon, Stock, Latest, Weekly, Monthly
$data = @(
  @("📉", "IIFLCAPS", -9.721500000000001, -5.9431, 20.4574),
  @("📉", "COHANCE", -8.5204, -10.2636, -10.5392),
  @("📉", "KHAICHEM", -8.4635, -6.7205, -4.6455),
  @("📉", "MOTILALOFS", -7.6514, -2.7408, 12.7798),
  @("📉", "VERANDA", -7.1035, -7.4096, 8.066700000000001),
  @("📉", "MOLDTKPAC", -7.0505, -1.2462, 0.466),
  @("📉", "NETWEB", -6.7177, 5.5715, 7.5565),
  @("📉", "FABTECH", -6.147, 15.6232, "N/A"),
  @("📉", "CARTRADE", -5.7253, 16.1193, 20.6359),
  @("📉", "TVSELECT", -5.6153, -0.9738, -2.9968),
  @("📉", "NSLNISP", -5.4542, 1.3037, 0.4681),
  @("📉", "NAM-INDIA", -5.1992, -7.0279, -1.2034),
  @("📉", "KICL", -5.0437, -0.8405, 21.7969),
  @("📉", "CCCL", -5.0146, -4.576, -12.4759),
  @("📉", "CREDITACC", -4.9692, -1.3216, 3.7319),
  @("📉", "KALAMANDIR", -4.8415, 1.7451, 25.9996),
  @("📉", "CRAMC", -4.7668, 5.978, "N/A"),
  @("📉", "SMLISUZU", -4.7654, 4.993, -2.8236),
  @("📉", "MANAKCOAT", -4.5883, -6.011, 24.8093),
  @("📉", "HDFCAMC", -4.401, -2.6247, -2.4311),
  @("📉", "ATHERENERG", -4.0945, -0.0142, 24.8806),
  @("📉", "SHAREINDIA", -4.0806, -1.6889, 54.7217),
  @("📉", "INDIQUBE", -4.0805, -4.8408, -6.7982),
  @("📉", "CHENNPETRO", -3.973, 4.8359, 6.8488),
  @("📉", "KFINTECH", -3.873, -1.7487, 7.3785),
  @("📉", "360ONE", -3.8488, -4.976, 10.0293),
  @("📉", "SMSPHARMA", -3.7339, -3.0871, 17.4387),
  @("📉", "BHARATWIRE", -3.5327, 22.8336, 23.8979),
  @("📉", "ABSLAMC", -3.5313, -5.9355, -1.2887),
  @("📉", "SUMMITSEC", -3.4113, -1.6476, 6.0097),
  @("📉", "SPLPETRO", -3.3984, -5.0241, -7.769),
  @("📉", "MPSLTD", -3.2782, -4.6314, 2.1755),
  @("📉", "CAMS", -3.2545, -0.6366000000000001, 2.5781),
  @("📉", "SPARC", -3.1709, 4.8337, 6.3311),
  @("📉", "PRUDENT", -3.127, -3.5103, 2.1213),
  @("📉", "ANANDRATHI", -3.0775, -0.8672, 9.1835),
  @("📉", "NLCINDIA", -3.0757, -4.5618, -11.6431),
  @("📉", "YATRA", -3.0403, -2.8455, 7.3711),
  @("📉", "BOSCHLTD", -3.0099, -3.123, -2.0055),
  @("📉", "DRREDDY", -2.9859, -2.5475, 2.2228),
  @("📉", "ROSSTECH", -2.9778, 1.9028, -6.8057),
  @("📉", "OAL", -2.9496, -1.278, 8.7362),
  @("📉", "ENDURANCE", -2.939, -2.2945, 3.4531),
  @("📉", "POLICYBZR", -2.907, 2.2365, 1.2573),
  @("📉", "AYMSYNTEX", -2.9052, -0.3705, -10.494),
  @("📉", "DIGITIDE", -2.8795, 3.2317, 6.2968),
  @("📉", "RUBICON", -2.8712, -2.9654, "N/A"),
  @("📉", "STARHEALTH", -2.8707, -1.5572, 7.5434),
  @("📉", "KIRIINDUS", -2.8411, -1.3849, 1.4335),
  @("📉", "UNIMECH", -2.8008, -1.6104, -0.4585),
  @("📉", "TTKPRESTIG", -2.7438, 8.001200000000001, 9.650499999999999),
  @("📉", "PFOCUS", -2.7039, -2.6276, -1.2163),
  @("📉", "ALLDIGI", -2.6342, -0.2306, -5.3103),
  @("📉", "PRIVISCL", -2.6288, -2.1048, 19.7451),
  @("📉", "CANHLIFE", -2.6148, 5.2953, "N/A"),
  @("📉", "GKENERGY", -2.6122, -9.807700000000001, 23.2758),
  @("📉", "SGFIN", -2.592, -0.06270000000000001, 11.7235),
  @("📉", "ARVINDFASN", -2.549, -2.9892, -4.4223),
  @("📉", "EDELWEISS", -2.5422, -3.3745, 8.5305),
  @("📉", "SAMHI", -2.5284, 1.8231, 2.8516),
  @("📉", "TBOTEK", -2.524, -3.5732, 1.036),
  @("📉", "UJJIVANSFB", -2.5201, 0.3845, 12.6645),
  @("📉", "AMBER", -2.5098, -0.1082, 2.763),
  @("📉", "GRPLTD", -2.4898, -5.9894, -5.4586),
  @("📉", "NESCO", -2.4722, 1.9934, 3.8931),
  @("📉", "PILANIINVS", -2.4546, -0.7907, 4.267),
  @("📉", "NSIL", -2.4088, -1.7646, 4.7431),
  @("📉", "COALINDIA", -2.4016, -3.058, -2.0387),
  @("📉", "JNKINDIA", -2.3482, -2.8371, 4.2622),
  @("📉", "FCL", -2.3453, -2.616, -0.02),
  @("📉", "DEEDEV", -2.3334, -6.6528, -7.4227),
  @("📉", "WEALTH", -2.2793, -3.8356, -2.7981),
  @("📉", "RATNAMANI", -2.2788, -0.4626, 0.8712),
  @("📉", "CSBBANK", -2.2695, 2.3137, 10.6999),
  @("📉", "BBOX", -2.2639, -4.7636, 5.1528)
)

$r = 2
foreach ($row in $data) {
    $losersWs.Cells.Item($r, 1).Value = $row[0]
    $losersWs.Cells.Item($r, 2).Value = $row[1]
    $losersWs.Cells.Item($r, 3).Value = $row[2]
    $losersWs.Cells.Item($r, 4).Value = $row[3]
    $losersWs.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Restore the originally active sheet selection
$metaWs.Select()
